$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wish list")
Write-Host "Before insert row131 D:" $ws.Cells.Item(131,4).Value()
$ws.Rows.Item(131).Insert()
Write-Host "After insert row131 D (should be blank now):" $ws.Cells.Item(131,4).Value()
Write-Host "After insert row132 D (should be old row131 val):" $ws.Cells.Item(132,4).Value()
